$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 25
# from serial date 45235 (2023-11-05) to 45236 (2023-11-06).
for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45235) {
        $cell.Value = 45236
    }
}
